# Fill in the two evidence sheets (B1, B2) that previously only held the
# generic placeholder text with the real tx-hash values, then leave the
# workbook focused on the newly completed "B1" sheet (matching the author's
# final view state) instead of the "Info" sheet.

$wb = $excel.ActiveWorkbook

# --- B1 sheet: first Interchain NFT-Transfer TxHash / Internal Transfer TxHash on IRISnet
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "33E21EEC017AC7E419A5946DC21D702BCF08C537D18746A142DF89E22FFA8FA5"
$wsB1.Range("A3").Value = "04395A7EDEE20A2CBB540C5E4AC8B78D57C1F9B2D73479712C01FDD781C44207"

# --- B2 sheet: same two pieces of evidence
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "A20F316898048289354F90065775737BE1E6E323A3B78F994B99E8652FEF65A3"
$wsB2.Range("A3").Value = "DC506E1B6EDA4526EE0C93B6C13661464F872967319EDBBA25394E2503451C20"

# Move the active/selected tab from "Info" to "B1" and park the cursor on A4
# (just past the now-filled A1:A3 evidence block), mirroring the saved view
# state in the workbook after upload.
$wsB1.Select() | Out-Null
$wsB1.Range("A4").Select() | Out-Null
